$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = -6
$ws.Range("F6").Value = 3
$ws.Range("F9").Value = -10
$ws.Range("F14").Value = -5
$ws.Range("F17").Value = -6
$ws.Range("F23").Value = -4
$ws.Range("F26").Value = 6
$ws.Range("F27").Value = -6
$ws.Range("F29").Value = -4
$ws.Range("F30").Value = 9
$ws.Range("F32").Value = 2
